$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Affiliations" label to lowercase "affiliations" (cleanup for create_creators kwargs)
$ws.Range("C7").Value = "affiliations"
$ws.Range("C8").Value = "affiliations"

# Update the active cell selection shown in the sheet view
$ws.Range("I16").Select()
